# Auto-generated: refresh market-price columns (H-N) for the affected Leve rows
# across sheets ALC, ARM, BSM, CRP, GSM, LTW, WVR, matching the scheduled-runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3705.7273
$ws.Range("I74").Value = 5003
$ws.Range("J74").Value = 3576
$ws.Range("K74").Value = 5003
$ws.Range("L74").Value = 3576
$ws.Range("M74").Value = -4067
$ws.Range("N74").Value = -5448
# Row 77
$ws.Range("H77").Value = 3705.7273
$ws.Range("I77").Value = 5003
$ws.Range("J77").Value = 3576
$ws.Range("K77").Value = 25015
$ws.Range("L77").Value = 17880
$ws.Range("M77").Value = -20335
$ws.Range("N77").Value = -27240
# Row 86
$ws.Range("H86").Value = 18535.738
$ws.Range("I86").Value = 1094.1538
$ws.Range("J86").Value = 41209.8
$ws.Range("K86").Value = 1094.1538
$ws.Range("L86").Value = 41209.8
$ws.Range("M86").Value = 28.84619999999995
$ws.Range("N86").Value = -43455.8
# Row 89
$ws.Range("H89").Value = 18535.738
$ws.Range("I89").Value = 1094.1538
$ws.Range("J89").Value = 41209.8
$ws.Range("K89").Value = 5470.769
$ws.Range("L89").Value = 206049
$ws.Range("M89").Value = 145.2309999999998
$ws.Range("N89").Value = -217281
# Row 113
$ws.Range("H113").Value = 6703.6562
$ws.Range("I113").Value = 2339.9375
$ws.Range("J113").Value = 11067.375
$ws.Range("K113").Value = 2339.9375
$ws.Range("L113").Value = 11067.375
$ws.Range("M113").Value = 914.0625
$ws.Range("N113").Value = -17575.375
# Row 132
$ws.Range("H132").Value = 19441.215
$ws.Range("J132").Value = 2183.2
$ws.Range("L132").Value = 6549.599999999999
$ws.Range("N132").Value = -11609.6

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3133.2222
$ws.Range("I45").Value = 3801.5
$ws.Range("J45").Value = 1796.6666
$ws.Range("K45").Value = 3801.5
$ws.Range("L45").Value = 1796.6666
$ws.Range("M45").Value = -3424.5
$ws.Range("N45").Value = -2550.6666
# Row 61
$ws.Range("H61").Value = 3147.963
$ws.Range("J61").Value = 2180
$ws.Range("L61").Value = 2180
$ws.Range("N61").Value = -2604
# Row 74
$ws.Range("H74").Value = 2981.9834
$ws.Range("I74").Value = 2956.9814
$ws.Range("J74").Value = 3207
$ws.Range("K74").Value = 2956.9814
$ws.Range("L74").Value = 3207
$ws.Range("M74").Value = -2082.9814
$ws.Range("N74").Value = -4955
# Row 77
$ws.Range("H77").Value = 2981.9834
$ws.Range("I77").Value = 2956.9814
$ws.Range("J77").Value = 3207
$ws.Range("K77").Value = 14784.907
$ws.Range("L77").Value = 16035
$ws.Range("M77").Value = -10416.907
$ws.Range("N77").Value = -24771
# Row 97
$ws.Range("H97").Value = 977.9091
$ws.Range("I97").Value = 976.8421
$ws.Range("K97").Value = 976.8421
$ws.Range("M97").Value = -480.8421
# Row 102
$ws.Range("H102").Value = 1232.8235
$ws.Range("I102").Value = 1232.8235
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1232.8235
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 389.1765
$ws.Range("N102").ClearContents()
# Row 136
$ws.Range("H136").Value = 3147.963
$ws.Range("J136").Value = 2180
$ws.Range("L136").Value = 6540
$ws.Range("N136").Value = -11640

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1617
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -674
$ws.Range("N85").ClearContents()
# Row 86
$ws.Range("H86").Value = 2238.3333
$ws.Range("I86").Value = 2569.0557
$ws.Range("K86").Value = 2569.0557
$ws.Range("M86").Value = -1446.0557
# Row 89
$ws.Range("H89").Value = 2238.3333
$ws.Range("I89").Value = 2569.0557
$ws.Range("K89").Value = 12845.2785
$ws.Range("M89").Value = -7229.2785
# Row 97
$ws.Range("H97").Value = 11832
$ws.Range("I97").Value = 7109.3335
$ws.Range("J97").Value = 26000
$ws.Range("K97").Value = 7109.3335
$ws.Range("L97").Value = 26000
$ws.Range("M97").Value = -6118.3335
$ws.Range("N97").Value = -27982
# Row 105
$ws.Range("H105").Value = 2582.1
$ws.Range("I105").Value = 2720.1667
$ws.Range("J105").Value = 2375
$ws.Range("K105").Value = 2720.1667
$ws.Range("L105").Value = 2375
$ws.Range("M105").Value = -973.1667000000002
$ws.Range("N105").Value = -5869

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2184.1167
$ws.Range("I31").Value = 1203.2222
$ws.Range("J31").Value = 3655.4583
$ws.Range("K31").Value = 1203.2222
$ws.Range("L31").Value = 3655.4583
$ws.Range("M31").Value = -908.2221999999999
$ws.Range("N31").Value = -4245.4583
# Row 34
$ws.Range("H34").Value = 2184.1167
$ws.Range("I34").Value = 1203.2222
$ws.Range("J34").Value = 3655.4583
$ws.Range("K34").Value = 1203.2222
$ws.Range("L34").Value = 3655.4583
$ws.Range("M34").Value = -1001.2222
$ws.Range("N34").Value = -4059.4583
# Row 105
$ws.Range("H105").Value = 1495.963
$ws.Range("I105").Value = 1790
$ws.Range("J105").Value = 996.1
$ws.Range("K105").Value = 1790
$ws.Range("L105").Value = 996.1
$ws.Range("M105").Value = -43
$ws.Range("N105").Value = -4490.1

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4981.365
$ws.Range("I70").Value = 4423.433
$ws.Range("J70").Value = 5742.1816
$ws.Range("K70").Value = 4423.433
$ws.Range("L70").Value = 5742.1816
$ws.Range("M70").Value = -4153.433
$ws.Range("N70").Value = -6282.1816
# Row 73
$ws.Range("H73").Value = 4981.365
$ws.Range("I73").Value = 4423.433
$ws.Range("J73").Value = 5742.1816
$ws.Range("K73").Value = 4423.433
$ws.Range("L73").Value = 5742.1816
$ws.Range("M73").Value = -3487.433
$ws.Range("N73").Value = -7614.1816
# Row 80
$ws.Range("H80").Value = 1904.826
$ws.Range("I80").Value = 1846.25
$ws.Range("J80").Value = 1968.7273
$ws.Range("K80").Value = 1846.25
$ws.Range("L80").Value = 1968.7273
$ws.Range("M80").Value = -848.25
$ws.Range("N80").Value = -3964.7273
# Row 83
$ws.Range("H83").Value = 1904.826
$ws.Range("I83").Value = 1846.25
$ws.Range("J83").Value = 1968.7273
$ws.Range("K83").Value = 9231.25
$ws.Range("L83").Value = 9843.636500000001
$ws.Range("M83").Value = -4239.25
$ws.Range("N83").Value = -19827.6365

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1032.3334
$ws.Range("J22").Value = 1822.5
$ws.Range("L22").Value = 1822.5
$ws.Range("N22").Value = -2412.5
# Row 27
$ws.Range("H27").Value = 1032.3334
$ws.Range("J27").Value = 1822.5
$ws.Range("L27").Value = 1822.5
$ws.Range("N27").Value = -2036.5
# Row 82
$ws.Range("H82").Value = 1576.8889
$ws.Range("J82").Value = 2010.5555
$ws.Range("L82").Value = 2010.5555
$ws.Range("N82").Value = -2732.5555
# Row 85
$ws.Range("H85").Value = 1576.8889
$ws.Range("J85").Value = 2010.5555
$ws.Range("L85").Value = 2010.5555
$ws.Range("N85").Value = -4506.5555
# Row 132
$ws.Range("H132").Value = 7767.8037
$ws.Range("I132").Value = 9414.517
$ws.Range("J132").Value = 5215.4
$ws.Range("K132").Value = 28243.551
$ws.Range("L132").Value = 15646.2
$ws.Range("M132").Value = -25713.551
$ws.Range("N132").Value = -20706.2
# Row 136
$ws.Range("H136").Value = 7247594.5
$ws.Range("I136").Value = 1259.091
$ws.Range("K136").Value = 3777.273
$ws.Range("M136").Value = -1227.273

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1181.6666
$ws.Range("I96").Value = 1022.5
$ws.Range("K96").Value = 1022.5
$ws.Range("M96").Value = 350.5
# Row 100
$ws.Range("H100").Value = 9574.583000000001
$ws.Range("I100").Value = 406
$ws.Range("J100").Value = 22410.6
$ws.Range("K100").Value = 812
$ws.Range("L100").Value = 44821.2
$ws.Range("M100").Value = -271
$ws.Range("N100").Value = -45903.2
# Row 126
$ws.Range("H126").Value = 49289.758
$ws.Range("I126").Value = 54777.04
$ws.Range("J126").Value = 1733.3334
$ws.Range("K126").Value = 164331.12
$ws.Range("L126").Value = 5200.0002
$ws.Range("M126").Value = -161861.12
$ws.Range("N126").Value = -10140.0002
# Row 132
$ws.Range("H132").Value = 1758.0392
$ws.Range("J132").Value = 2162.077
$ws.Range("L132").Value = 6486.231000000001
$ws.Range("N132").Value = -11546.231
# Row 136
$ws.Range("H136").Value = 4904677
$ws.Range("I136").Value = 7246824.5
$ws.Range("K136").Value = 21740473.5
$ws.Range("M136").Value = -21737923.5
